# Apply the "Add 2 more test case for BOA Homepage" change.
$wb = $excel.ActiveWorkbook

# Rename the BMWUSA sheet to countyList
$ws = $wb.Worksheets.Item("BMWUSA")
$ws.Name = "countyList"

# Populate the countyList sheet with County / UNICOI, TN / WASHINGTON, TN / CARTER, TN
$ws.Range("A1").Value = "County"
$ws.Range("A2").Value = "UNICOI, TN"
$ws.Range("A3").Value = "WASHINGTON, TN"
$ws.Range("A4").Value = "CARTER, TN"

# Match style used for the interior rows (A2:A3) -- reuse the existing
# "code" cell style (Menlo font) already present on NavigationBarMenu!A2
# by copying formats, instead of building a brand-new style entry.
$nav = $wb.Worksheets.Item("NavigationBarMenu")
$nav.Range("A2").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Column width
$ws.Columns.Item(1).ColumnWidth = 16

# Make this sheet the active/selected tab, then select A4
$ws.Activate() | Out-Null
$ws.Range("A4").Select() | Out-Null
